# Auto-generated edit script: updates currentAveragePrice / Leve price / profit
# columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve tables to
# reflect refreshed market-board averages from the scheduled data-pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 549.5333000000001
$ws.Range("I18").Value = 495.9
$ws.Range("J18").Value = 656.8
$ws.Range("K18").Value = 495.9
$ws.Range("L18").Value = 656.8
$ws.Range("M18").Value = -211.9
$ws.Range("N18").Value = -1224.8
$ws.Range("H40").Value = 1215.3846
$ws.Range("I40").Value = 879.4
$ws.Range("J40").Value = 1425.375
$ws.Range("K40").Value = 879.4
$ws.Range("L40").Value = 1425.375
$ws.Range("M40").Value = -704.4
$ws.Range("N40").Value = -1775.375
$ws.Range("H61").Value = 177.14285
$ws.Range("I61").Value = 177.14285
$ws.Range("K61").Value = 531.4285500000001
$ws.Range("M61").Value = -359.4285500000001
$ws.Range("H74").Value = 6333.0835
$ws.Range("I74").Value = 5119.8
$ws.Range("J74").Value = 7199.7144
$ws.Range("K74").Value = 5119.8
$ws.Range("L74").Value = 7199.7144
$ws.Range("M74").Value = -4183.8
$ws.Range("N74").Value = -9071.714400000001
$ws.Range("H77").Value = 6333.0835
$ws.Range("I77").Value = 5119.8
$ws.Range("J77").Value = 7199.7144
$ws.Range("K77").Value = 25599
$ws.Range("L77").Value = 35998.572
$ws.Range("M77").Value = -20919
$ws.Range("N77").Value = -45358.572
$ws.Range("H137").Value = 1061923.4
$ws.Range("I137").Value = 3406708.2
$ws.Range("J137").Value = 2988.2258
$ws.Range("K137").Value = 10220124.6
$ws.Range("L137").Value = 8964.6774
$ws.Range("M137").Value = -10217574.6
$ws.Range("N137").Value = -14064.6774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3363.6924
$ws.Range("I45").Value = 2525.4443
$ws.Range("K45").Value = 2525.4443
$ws.Range("M45").Value = -2148.4443
$ws.Range("H61").Value = 1608.75
$ws.Range("I61").Value = 1512.5
$ws.Range("J61").Value = 1833.3334
$ws.Range("K61").Value = 1512.5
$ws.Range("L61").Value = 1833.3334
$ws.Range("M61").Value = -1300.5
$ws.Range("N61").Value = -2257.3334
$ws.Range("H136").Value = 1608.75
$ws.Range("I136").Value = 1512.5
$ws.Range("J136").Value = 1833.3334
$ws.Range("K136").Value = 4537.5
$ws.Range("L136").Value = 5500.0002
$ws.Range("M136").Value = -1987.5
$ws.Range("N136").Value = -10600.0002
$ws.Range("H137").Value = 43610
$ws.Range("J137").Value = 43610
$ws.Range("L137").Value = 43610
$ws.Range("N137").Value = -53810

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 159.91667
$ws.Range("I22").Value = 156.27272
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 156.27272
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 16.72728000000001
$ws.Range("N22").Value = -546
$ws.Range("H86").Value = 2638.1428
$ws.Range("I86").Value = 1980
$ws.Range("J86").Value = 2747.8333
$ws.Range("K86").Value = 1980
$ws.Range("L86").Value = 2747.8333
$ws.Range("M86").Value = -857
$ws.Range("N86").Value = -4993.8333
$ws.Range("H89").Value = 2638.1428
$ws.Range("I89").Value = 1980
$ws.Range("J89").Value = 2747.8333
$ws.Range("K89").Value = 9900
$ws.Range("L89").Value = 13739.1665
$ws.Range("M89").Value = -4284
$ws.Range("N89").Value = -24971.1665
$ws.Range("H94").Value = 664
$ws.Range("J94").Value = 1100
$ws.Range("L94").Value = 1100
$ws.Range("N94").Value = -2002
$ws.Range("H99").Value = 1276
$ws.Range("I99").Value = 1188.8
$ws.Range("J99").Value = 1566.6666
$ws.Range("K99").Value = 1188.8
$ws.Range("L99").Value = 1566.6666
$ws.Range("M99").Value = 309.2
$ws.Range("N99").Value = -4562.6666
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H134").Value = 2667.756
$ws.Range("I134").Value = 1095.85
$ws.Range("J134").Value = 4164.8096
$ws.Range("K134").Value = 3287.55
$ws.Range("L134").Value = 12494.4288
$ws.Range("M134").Value = -752.5499999999997
$ws.Range("N134").Value = -17564.4288

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3187.121
$ws.Range("I132").Value = 2850.7
$ws.Range("J132").Value = 3704.6924
$ws.Range("K132").Value = 8552.099999999999
$ws.Range("L132").Value = 11114.0772
$ws.Range("M132").Value = -6022.099999999999
$ws.Range("N132").Value = -16174.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 456.32074
$ws.Range("I113").Value = 452.8857
$ws.Range("K113").Value = 1358.6571
$ws.Range("M113").Value = 811.3429000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 24755.084
$ws.Range("J46").Value = 25642.363
$ws.Range("L46").Value = 25642.363
$ws.Range("N46").Value = -25954.363
$ws.Range("H70").Value = 6912.1562
$ws.Range("I70").Value = 6195.2
$ws.Range("K70").Value = 6195.2
$ws.Range("M70").Value = -5925.2
$ws.Range("H73").Value = 6912.1562
$ws.Range("I73").Value = 6195.2
$ws.Range("K73").Value = 6195.2
$ws.Range("M73").Value = -5259.2
$ws.Range("H132").Value = 4529.1143
$ws.Range("I132").Value = 3894.5625
$ws.Range("J132").Value = 5063.4736
$ws.Range("K132").Value = 11683.6875
$ws.Range("L132").Value = 15190.4208
$ws.Range("M132").Value = -9153.6875
$ws.Range("N132").Value = -20250.4208
$ws.Range("H137").Value = 39086.668
$ws.Range("J137").Value = 48630
$ws.Range("L137").Value = 48630
$ws.Range("N137").Value = -58830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1630.2413
$ws.Range("I46").Value = 991.93335
$ws.Range("J46").Value = 2314.1428
$ws.Range("K46").Value = 991.93335
$ws.Range("L46").Value = 2314.1428
$ws.Range("M46").Value = -803.93335
$ws.Range("N46").Value = -2690.1428
$ws.Range("H93").Value = 11113950
$ws.Range("I93").Value = 13891187
$ws.Range("K93").Value = 13891187
$ws.Range("M93").Value = -13889939
$ws.Range("H136").Value = 4954.16
$ws.Range("J136").Value = 8724.091
$ws.Range("L136").Value = 26172.273
$ws.Range("N136").Value = -31272.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H132").Value = 9806651
$ws.Range("I132").Value = 1047.125
$ws.Range("J132").Value = 12823760
$ws.Range("K132").Value = 3141.375
$ws.Range("L132").Value = 38471280
$ws.Range("M132").Value = -611.375
$ws.Range("N132").Value = -38476340
